$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix typo in header: "Акутальность" -> "Актуальность"
$ws.Range("H1").Value = "Актуальность"

# Fill in the "Пласт" (formation) column G for data rows with "Юг1"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 7).Value = "Юг1"
}
